# Refatorando artefatos de acordo com o feedback da ac4
#
# This script normalises the bullet list terminators in the "Usuários e
# Outros Stakeholders" table: two bullets that were missing a closing
# full stop gain one (as a brand-new trailing run), and three bullets
# that incorrectly ended in a comma have that comma swapped for a full
# stop (keeping it as its own trailing run).
#
# Because any in-place text mutation on a paragraph causes this host's
# Word engine to re-flow that paragraph's runs into a single run, we
# nudge the edited character range through a harmless Bold on/off
# toggle immediately afterwards. That forces the engine to keep (or
# recreate) a distinct trailing <w:r> for the punctuation, mirroring
# the two-run shape the diff expects, without leaving any visible
# formatting behind.

$d = $word.ActiveDocument

function Add-TrailingPeriod($anchorText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find anchor text: $anchorText"
    }
    # Insert a brand-new run containing "." right after the matched text.
    $insertionPoint = $d.Range($rng.End, $rng.End)
    $insertionPoint.InsertAfter(".")

    # Re-touch just the inserted character so it is kept as its own run
    # instead of being folded back into the preceding one.
    $newRun = $d.Range($rng.End, $rng.End + 1)
    $newRun.Bold = 1
    $newRun.Bold = 0
}

function Replace-TrailingCommaWithPeriod($anchorText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find anchor text: $anchorText"
    }
    # The match ends right after the comma; replace just that one
    # character in place.
    $commaRng = $d.Range($rng.End - 1, $rng.End)
    if ($commaRng.Text -ne ",") {
        throw "Expected a trailing comma near: $anchorText"
    }
    $commaRng.Text = "."

    # Re-touch the replaced character so it stays split out as its own
    # trailing run, matching the original document's run layout.
    $newRun = $d.Range($rng.End - 1, $rng.End)
    $newRun.Bold = 1
    $newRun.Bold = 0
}

# 1) "Verificar os orçamentos recebidos, enviados, faturados e afins" -> add "."
Add-TrailingPeriod("Verificar os orçamentos recebidos, enviados, faturados e afins")

# 2) "Verificar os pedidos" -> add "."
Add-TrailingPeriod("Verificar os pedidos")

# 3) "Acompanhar os processos de produção," -> "...produção."
Replace-TrailingCommaWithPeriod("Acompanhar os processos de produção,")

# 4) "...ao produto," -> "...ao produto."
Replace-TrailingCommaWithPeriod(" ao produto,")

# 5) "...orçamento solicitado," -> "...orçamento solicitado."
Replace-TrailingCommaWithPeriod("orçamento solicitado,")
